$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("relative_to_now") description: clarify that large initial imports
# can trigger HTTPErrors and that the time interval should be reduced in
# that case. This is the text in C7.
$ws.Range("C7").Value = "If FALSE, import data in time interval <from_time> to <to_time>. Use for initial import of large data sets (if getting HTTPErrors, reduce the time intervall). If TRUE, import the last x hours, where x = <relative offset>. Use this option for continuous updates."

# The longer text now needs more vertical space to display fully wrapped.
$ws.Rows.Item(7).RowHeight = 105

# The last selected/active cell moved to C8.
$ws.Range("C8").Select()
